$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:category-name" (column E) and "codeforiati:group-code"
# (column F) values were swapped/mislabeled in the source data (including
# the header row). Fix this by swapping the content of columns E and F for
# every row in the used range.
#
# All of these cells hold text (shared strings) in the source file -
# including values that look numeric, like the group codes ("110", "120",
# ...). Writing such a numeric-looking string back through .Value2 would
# turn it into a real number, so numeric-looking strings are written back
# through .Formula with a leading apostrophe to force a text cell, exactly
# like typing '110 into Excel.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

function Get-CellText($cell) {
    $v = $cell.Value2
    return [string]$v
}

function Set-CellText($cell, $text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cell.Formula = "'" + $text
    } else {
        $cell.Value2 = $text
    }
}

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eText = Get-CellText $eCell
    $fText = Get-CellText $fCell

    Set-CellText $eCell $fText
    Set-CellText $fCell $eText
}
